# Replenishment_Requests.xlsx - add a new pending replenishment request row
# (Request 4, pharmacist P001, PARACETAMOL, qty 10, status PENDING, dated 02/11/2024)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple columns - Excel will pick up the existing column formatting automatically.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "P001"
$ws.Range("D5").Value = "PARACETAMOL"
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = "PENDING"

# The date column needs special handling: assigning the literal string "02/11/2024"
# directly would make Excel auto-recognise it as a date value (and pick up a date
# number format), but the source data keeps it as plain text. Entering it as a
# text-formula and then converting that formula to its value keeps the cell as
# plain text using the existing (unformatted) column style instead of creating a
# new date-formatted style.
$ws.Range("C5").Formula = '="02/11/2024"'
$ws.Range("C5").Copy()
$ws.Range("C5").PasteSpecial(-4163)
